$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (2023-10-05) for rows 2-70
# that must be updated to the new date (2023-10-08), keeping existing formatting.
for ($row = 2; $row -le 70; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45207
    }
}
